# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# "Bad Drivers" section - Intel(R) Wi-Fi 6 AX201 160MHz - 22.150.0.3
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 98.90000000000001

# Totals row
$ws.Range("C4").Value = 6

# "Good Drivers" section - Total Samples updates
$ws.Range("B12").Value = 449371
$ws.Range("B13").Value = 77999
